# Weekly refresh of the "Fruta, Vega Central Mapocho de Santiago - Mora" data.
# The 11 daily records (rows 2-12) are re-shuffled/updated in place: each row's
# Fecha (D), Volumen (M), Precio minimo/maximo/promedio (N/O/P), Origen (R) and
# Precio $/Kg (S) are updated to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó";  S = 1600 },
    @{ Row = 3;  D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 },
    @{ Row = 4;  D = 44232; M = 200; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó";  S = 1500 },
    @{ Row = 5;  D = 44188; M = 150; N = 3000; O = 3400; P = 3240; R = "Provincia de Linares"; S = 1620 },
    @{ Row = 6;  D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = "Provincia de Curicó";  S = 1700 },
    @{ Row = 7;  D = 44237; M = 100; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 },
    @{ Row = 8;  D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 },
    @{ Row = 9;  D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 },
    @{ Row = 10; D = 44194; M = 120; N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 },
    @{ Row = 11; D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 },
    @{ Row = 12; D = 44533; M = 150; N = 4000; O = 4000; P = 4000; R = "Provincia de Curicó";  S = 2000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $item.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $item.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $item.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $item.S   # S: Precio $/Kg
}
